$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the value of B7 (previously "運動鞋") while keeping its existing
# cell formatting/style - mirrors pressing Delete on a selected cell.
$ws.Range("B7").ClearContents()

# Move the active cell / selection to C13 (matches the saved sheet view).
$ws.Range("C13").Select()
